$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.04983966666666667
$ws.Range("H2").Value = 0.149519
$ws.Range("I2").Value = 0.1823731600337622
$ws.Range("J2").Value = 0.1823731600337622
$ws.Range("Q2").Value = 0.0003533964631111111
$ws.Range("R2").Value = 0.003180568168
$ws.Range("S2").Value = 0.1823731600337622
$ws.Range("T2").Value = 0.1823731600337622

# Row 3 updates
$ws.Range("I3").Value = 0.4031422744592926
$ws.Range("J3").Value = 0.4031422744592926
$ws.Range("S3").Value = 0.4031422744592926
$ws.Range("T3").Value = 0.4031422744592926

# Row 4 updates
$ws.Range("I4").Value = 0.4144845655069452
$ws.Range("J4").Value = 0.4144845655069451
$ws.Range("S4").Value = 0.4144845655069452
$ws.Range("T4").Value = 0.4144845655069451
